$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.7065936666666666
$ws.Range("H2").Value = 2.119781
$ws.Range("I2").Value = 0.005187843618793344
$ws.Range("J2").Value = 0.005187843618793344
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.827738000000002
$ws.Range("N2").Value = 29.483214
$ws.Range("O2").Value = 0.3869625527756497
$ws.Range("P2").Value = 0.3869625527756497
$ws.Range("Q2").Value = 6.944217428459334
$ws.Range("R2").Value = 62.497956856134
$ws.Range("S2").Value = 0.002007501210129137
$ws.Range("T2").Value = 0.002007501210129137

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.7065936666666666
$ws.Range("H3").Value = 2.119781
$ws.Range("I3").Value = 0.005187843618793344
$ws.Range("J3").Value = 0.005187843618793344
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.362736333333332
$ws.Range("N3").Value = 28.088209
$ws.Range("O3").Value = 0.368653331266258
$ws.Range("P3").Value = 0.368653331266258
$ws.Range("Q3").Value = 6.61565019580322
$ws.Range("R3").Value = 59.54085176222899
$ws.Range("S3").Value = 0.001912515832156565
$ws.Range("T3").Value = 0.001912515832156565

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.7065936666666666
$ws.Range("H4").Value = 2.119781
$ws.Range("I4").Value = 0.005187843618793344
$ws.Range("J4").Value = 0.005187843618793344
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.206655
$ws.Range("N4").Value = 18.619965
$ws.Range("O4").Value = 0.2443841159580923
$ws.Range("P4").Value = 0.2443841159580923
$ws.Range("Q4").Value = 4.385583114185
$ws.Range("R4").Value = 39.47024802766499
$ws.Range("S4").Value = 0.001267826576507641
$ws.Range("T4").Value = 0.001267826576507641

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 111.9320066666667
$ws.Range("H5").Value = 335.79602
$ws.Range("I5").Value = 0.8218100075305903
$ws.Range("J5").Value = 0.8218100075305903
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.827738000000002
$ws.Range("N5").Value = 29.483214
$ws.Range("O5").Value = 0.3869625527756497
$ws.Range("P5").Value = 0.3869625527756497
$ws.Range("Q5").Value = 1100.038435334254
$ws.Range("R5").Value = 9900.345918008281
$ws.Range("S5").Value = 0.3180096984106131
$ws.Range("T5").Value = 0.3180096984106131

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 111.9320066666667
$ws.Range("H6").Value = 335.79602
$ws.Range("I6").Value = 0.8218100075305903
$ws.Range("J6").Value = 0.8218100075305903
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.362736333333332
$ws.Range("N6").Value = 28.088209
$ws.Range("O6").Value = 0.368653331266258
$ws.Range("P6").Value = 0.368653331266258
$ws.Range("Q6").Value = 1047.989865680909
$ws.Range("R6").Value = 9431.90879112818
$ws.Range("S6").Value = 0.3029629969441007
$ws.Range("T6").Value = 0.3029629969441007

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 111.9320066666667
$ws.Range("H7").Value = 335.79602
$ws.Range("I7").Value = 0.8218100075305903
$ws.Range("J7").Value = 0.8218100075305903
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.206655
$ws.Range("N7").Value = 18.619965
$ws.Range("O7").Value = 0.2443841159580923
$ws.Range("P7").Value = 0.2443841159580923
$ws.Range("Q7").Value = 694.7233488377001
$ws.Range("R7").Value = 6252.5101395393
$ws.Range("S7").Value = 0.2008373121758764
$ws.Range("T7").Value = 0.2008373121758764

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.563205
$ws.Range("H8").Value = 70.689615
$ws.Range("I8").Value = 0.1730021488506163
$ws.Range("J8").Value = 0.1730021488506163
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.827738000000002
$ws.Range("N8").Value = 29.483214
$ws.Range("O8").Value = 0.3869625527756497
$ws.Range("P8").Value = 0.3869625527756497
$ws.Range("Q8").Value = 231.57300518029
$ws.Range("R8").Value = 2084.15704662261
$ws.Range("S8").Value = 0.06694535315490741
$ws.Range("T8").Value = 0.06694535315490742

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.563205
$ws.Range("H9").Value = 70.689615
$ws.Range("I9").Value = 0.1730021488506163
$ws.Range("J9").Value = 0.1730021488506163
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.362736333333332
$ws.Range("N9").Value = 28.088209
$ws.Range("O9").Value = 0.368653331266258
$ws.Range("P9").Value = 0.368653331266258
$ws.Range("Q9").Value = 220.6160755832817
$ws.Range("R9").Value = 1985.544680249535
$ws.Range("S9").Value = 0.06377781849000072
$ws.Range("T9").Value = 0.06377781849000075

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.563205
$ws.Range("H10").Value = 70.689615
$ws.Range("I10").Value = 0.1730021488506163
$ws.Range("J10").Value = 0.1730021488506163
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.206655
$ws.Range("N10").Value = 18.619965
$ws.Range("O10").Value = 0.2443841159580923
$ws.Range("P10").Value = 0.2443841159580923
$ws.Range("Q10").Value = 146.248684129275
$ws.Range("R10").Value = 1316.238157163475
$ws.Range("S10").Value = 0.04227897720570815
$ws.Range("T10").Value = 0.04227897720570815
